$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.137.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.10"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3884"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08236"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.180"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.16"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.161"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.81"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06664"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.46"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.899"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.125.75"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.240"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.069.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.38"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = "  -5.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1041"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.020"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.760"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02405"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06402"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.049"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.69%  "
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.241"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.176"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6359"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.906"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5955"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.80"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.679"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.274"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.958"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.194"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06799"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.56"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.24%  "
